$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain text values (e.g. "1.005", "30.473.47").
# Mark the cells as Text before writing so Excel does not reinterpret
# numeric-looking strings as floating point numbers.

$ws.Range("D2").Value = "30.473.47"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "2.108.20"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.65"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5242"
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4585"
$ws.Range("E8").Value = "  +5.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.37"
$ws.Range("E9").Value = "  +12.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08962"
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.178"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.41"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "2.109.42"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.780"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.820"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.59"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001130"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06632"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.25"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.287"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "30.544.34"
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.34"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.354"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("D26").Value = "2.352.27"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.30"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.563"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.30"
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.78"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.194"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.692"
$ws.Range("E33").Value = "  +9.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.151"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.937"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.41"
$ws.Range("E36").Value = "  +8.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02570"
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06820"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.543"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.79"
$ws.Range("E40").Value = "  +1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2286"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6885"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.246"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.345"
$ws.Range("E44").Value = "  +5.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.003"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.04"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6372"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.653"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  +24.81%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.221"
$ws.Range("E51").Value = "  +2.28%  "
